$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CutfillerCoefficient")

# Locate the XML-mapped table (Tabela6) backing the CutfillerCoefficient sheet.
$lo = $ws.ListObjects.Item(1)

# Remove the (now unused) "Title" column - the CutfillerCoefficient.Title field
# was dropped from the schema, so the leading Title column goes away and the
# remaining CFTProductivityRateMin / CFTProductivityRateMax columns shift left.
$ws.Columns("A").Select() | Out-Null
$ws.Columns("A").Delete()

# Resync the table's range/column count with the shrunk data.
$lo.Resize($ws.Range("A1:B2"))

# Resync the list-column names from the (already-correct) header cells so the
# table definition reflects CFTProductivityRateMin / CFTProductivityRateMax
# instead of the stale Title / CFTProductivityRateMin names.
$lo.HeaderRowRange.Cells(1, 1).Value2 = "CFTProductivityRateMin"
$lo.HeaderRowRange.Cells(1, 2).Value2 = "CFTProductivityRateMax"

# Make CutfillerCoefficient the active sheet/tab (was Usage before).
$ws.Activate()
